$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3695477391476345
$ws.Cells.Item(2, 3).Value = 0.03713283776956189
$ws.Cells.Item(2, 5).Value = 0.1469818966848848
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.8271985361931655
$ws.Cells.Item(2, 8).Value = 0.8884377347848016
$ws.Cells.Item(2, 11).Value = 0.329930102009456
$ws.Cells.Item(2, 13).Value = 0.2427936730292899
$ws.Cells.Item(2, 14).Value = 1.7534874407595

$ws.Cells.Item(3, 2).Value = 0.3353331068085197
$ws.Cells.Item(3, 3).Value = 0.03282023116986466
$ws.Cells.Item(3, 5).Value = 0.136066826490385
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.8254898051503545
$ws.Cells.Item(3, 8).Value = 0.8923437923269688
$ws.Cells.Item(3, 11).Value = 0.2940628102415133
$ws.Cells.Item(3, 13).Value = 0.220595884493477
$ws.Cells.Item(3, 14).Value = 1.772785254169563

$ws.Cells.Item(4, 2).Value = 0.3144587076047003
$ws.Cells.Item(4, 3).Value = 0.03015878146882756
$ws.Cells.Item(4, 5).Value = 0.1294567659513888
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.8249969203230876
$ws.Cells.Item(4, 8).Value = 0.8951597997564562
$ws.Cells.Item(4, 11).Value = 0.2721185965840647
$ws.Cells.Item(4, 13).Value = 0.2070775817920847
$ws.Cells.Item(4, 14).Value = 1.785265768905443

$ws.Cells.Item(5, 2).Value = 0.3059859284374511
$ws.Cells.Item(5, 3).Value = 0.02907081380260479
$ws.Cells.Item(5, 5).Value = 0.1267859904782824
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.8249356101334513
$ws.Cells.Item(5, 8).Value = 0.8964123173796708
$ws.Cells.Item(5, 11).Value = 0.2631959604163825
$ws.Cells.Item(5, 13).Value = 0.2015965581517563
$ws.Cells.Item(5, 14).Value = 1.790510346791862

$ws.Cells.Item(6, 2).Value = 0.3045810697681475
$ws.Cells.Item(6, 3).Value = 0.02888995156958174
$ws.Cells.Item(6, 5).Value = 0.1263438857434807
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.8249338479196098
$ws.Cells.Item(6, 8).Value = 0.8966266354217822
$ws.Cells.Item(6, 11).Value = 0.2617155624131158
$ws.Cells.Item(6, 13).Value = 0.2006881120033484
$ws.Cells.Item(6, 14).Value = 1.79139078372792

$ws.Cells.Item(7, 2).Value = 0.3143443040817999
$ws.Cells.Item(7, 3).Value = 0.03014412256089827
$ws.Cells.Item(7, 5).Value = 0.1294206545909091
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.8249955289118844
$ws.Cells.Item(7, 8).Value = 0.8951762666969927
$ws.Cells.Item(7, 11).Value = 0.2719981824744053
$ws.Cells.Item(7, 13).Value = 0.2070035504810619
$ws.Cells.Item(7, 14).Value = 1.785335856883687

$ws.Cells.Item(8, 2).Value = 0.3577228890834192
$ws.Cells.Item(8, 3).Value = 0.03564864231852027
$ws.Cells.Item(8, 5).Value = 0.1431991866071343
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.8264937028314279
$ws.Cells.Item(8, 8).Value = 0.8896978208806132
$ws.Cells.Item(8, 11).Value = 0.3175468155848193
$ws.Cells.Item(8, 13).Value = 0.2351166772835782
$ws.Cells.Item(8, 14).Value = 1.760010154600643

$ws.Cells.Item(9, 2).Value = 0.4438464404831848
$ws.Cells.Item(9, 3).Value = 0.04633701988805683
$ws.Cells.Item(9, 5).Value = 0.1709587061968278
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.8338620999660549
$ws.Cells.Item(9, 8).Value = 0.8822717231495716
$ws.Cells.Item(9, 11).Value = 0.4074911172485258
$ws.Cells.Item(9, 13).Value = 0.2911398567359882
$ws.Cells.Item(9, 14).Value = 1.7153637696657

$ws.Cells.Item(10, 2).Value = 0.5077718763994312
$ws.Cells.Item(10, 3).Value = 0.05412740362463353
$ws.Cells.Item(10, 5).Value = 0.1918222537306278
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.8420016423054903
$ws.Cells.Item(10, 8).Value = 0.8788430688133957
$ws.Cells.Item(10, 11).Value = 0.4739618658029201
$ws.Cells.Item(10, 13).Value = 0.3328640751312264
$ws.Cells.Item(10, 14).Value = 1.685627825624394

$ws.Cells.Item(11, 2).Value = 0.5369960320650193
$ws.Cells.Item(11, 3).Value = 0.05765843565886541
$ws.Cells.Item(11, 5).Value = 0.201419302115859
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.8463019957392248
$ws.Cells.Item(11, 8).Value = 0.8777246789844213
$ws.Cells.Item(11, 11).Value = 0.5042877981450715
$ws.Cells.Item(11, 13).Value = 0.3519723954598888
$ws.Cells.Item(11, 14).Value = 1.672767423605485

$ws.Cells.Item(12, 2).Value = 0.5480831648260107
$ws.Cells.Item(12, 3).Value = 0.05899372367242961
$ws.Cells.Item(12, 5).Value = 0.2050689858380252
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.8480167815635298
$ws.Cells.Item(12, 8).Value = 0.8773647240050337
$ws.Cells.Item(12, 11).Value = 0.5157841321405101
$ws.Cells.Item(12, 13).Value = 0.3592268555748959
$ws.Cells.Item(12, 14).Value = 1.667993557797395

$ws.Cells.Item(13, 2).Value = 0.545694439517689
$ws.Cells.Item(13, 3).Value = 0.05870622703879746
$ws.Cells.Item(13, 5).Value = 0.2042822685540102
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.847643625642192
$ws.Cells.Item(13, 8).Value = 0.8774394186217762
$ws.Cells.Item(13, 11).Value = 0.5133076338270257
$ws.Cells.Item(13, 13).Value = 0.3576636496955885
$ws.Cells.Item(13, 14).Value = 1.669017417925826

$ws.Cells.Item(14, 2).Value = 0.5379077655658762
$ws.Cells.Item(14, 3).Value = 0.05776832736292192
$ws.Cells.Item(14, 5).Value = 0.2017192520298536
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.8464413396380195
$ws.Cells.Item(14, 8).Value = 0.877693790995707
$ws.Cells.Item(14, 11).Value = 0.5052333572807584
$ws.Cells.Item(14, 13).Value = 0.3525688513669323
$ws.Cells.Item(14, 14).Value = 1.672372747694002

$ws.Cells.Item(15, 2).Value = 0.533140878532862
$ws.Cells.Item(15, 3).Value = 0.05719359780667332
$ws.Cells.Item(15, 5).Value = 0.2001513544885754
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.8457161599360035
$ws.Cells.Item(15, 8).Value = 0.877857880950188
$ws.Cells.Item(15, 11).Value = 0.5002892626801554
$ws.Cells.Item(15, 13).Value = 0.349450561958875
$ws.Cells.Item(15, 14).Value = 1.674440504589796

$ws.Cells.Item(16, 2).Value = 0.5058648970035904
$ws.Cells.Item(16, 3).Value = 0.05389638468574276
$ws.Cells.Item(16, 5).Value = 0.191197219046046
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.8417326638618761
$ws.Cells.Item(16, 8).Value = 0.8789250466973613
$ws.Cells.Item(16, 11).Value = 0.4719817580654535
$ws.Cells.Item(16, 13).Value = 0.3316178878948435
$ws.Cells.Item(16, 14).Value = 1.686481716890714

$ws.Cells.Item(17, 2).Value = 0.4891687800806039
$ws.Cells.Item(17, 3).Value = 0.05187036711694759
$ws.Cells.Item(17, 5).Value = 0.1857314894292799
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.8394422680123341
$ws.Cells.Item(17, 8).Value = 0.8796928174912182
$ws.Cells.Item(17, 11).Value = 0.4546385018867625
$ws.Cells.Item(17, 13).Value = 0.3207109722304793
$ws.Cells.Item(17, 14).Value = 1.694039511795431

$ws.Cells.Item(18, 2).Value = 0.4795791925896253
$ws.Cells.Item(18, 3).Value = 0.05070384497643943
$ws.Cells.Item(18, 5).Value = 0.1825977206171174
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.8381811213188968
$ws.Cells.Item(18, 8).Value = 0.8801759525910171
$ws.Cells.Item(18, 11).Value = 0.4446714137723689
$ws.Cells.Item(18, 13).Value = 0.3144496176579992
$ws.Cells.Item(18, 14).Value = 1.698449287754659

$ws.Cells.Item(19, 2).Value = 0.4763346595512701
$ws.Cells.Item(19, 3).Value = 0.05030867240114389
$ws.Cells.Item(19, 5).Value = 0.1815383864334592
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.8377637639637641
$ws.Cells.Item(19, 8).Value = 0.8803466644619391
$ws.Cells.Item(19, 11).Value = 0.4412981566096619
$ws.Cells.Item(19, 13).Value = 0.3123316886297189
$ws.Cells.Item(19, 14).Value = 1.699953129538514

$ws.Cells.Item(20, 2).Value = 0.490944706113936
$ws.Cells.Item(20, 3).Value = 0.05208616534163468
$ws.Cells.Item(20, 5).Value = 0.1863122917112605
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.839680262262533
$ws.Cells.Item(20, 8).Value = 0.8796067877277238
$ws.Cells.Item(20, 11).Value = 0.4564838634383648
$ws.Cells.Item(20, 13).Value = 0.3218707872224655
$ws.Cells.Item(20, 14).Value = 1.693228477421471

$ws.Cells.Item(21, 2).Value = 0.5401943431312191
$ws.Cells.Item(21, 3).Value = 0.05804386100351167
$ws.Cells.Item(21, 5).Value = 0.2024716498811614
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.8467921339002231
$ws.Cells.Item(21, 8).Value = 0.877617350047359
$ws.Cells.Item(21, 11).Value = 0.5076046279036461
$ws.Cells.Item(21, 13).Value = 0.3540648121347729
$ws.Cells.Item(21, 14).Value = 1.671384596324874

$ws.Cells.Item(22, 2).Value = 0.5725017109093358
$ws.Cells.Item(22, 3).Value = 0.061926856321179
$ws.Cells.Item(22, 5).Value = 0.2131231186412208
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.8519435664186972
$ws.Cells.Item(22, 8).Value = 0.8766876289835892
$ws.Cells.Item(22, 11).Value = 0.5410881754395405
$ws.Cells.Item(22, 13).Value = 0.3752137000382092
$ws.Cells.Item(22, 14).Value = 1.657668548467985

$ws.Cells.Item(23, 2).Value = 0.5552477406201035
$ws.Cells.Item(23, 3).Value = 0.05985540427622027
$ws.Cells.Item(23, 5).Value = 0.2074298818044298
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.8491479562651989
$ws.Cells.Item(23, 8).Value = 0.8771499073744451
$ws.Cells.Item(23, 11).Value = 0.5232107135437616
$ws.Cells.Item(23, 13).Value = 0.3639161770237891
$ws.Cells.Item(23, 14).Value = 1.664937727144213

$ws.Cells.Item(24, 2).Value = 0.4901417811861677
$ws.Cells.Item(24, 3).Value = 0.05198860837867869
$ws.Cells.Item(24, 5).Value = 0.1860496843822972
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.8395724917961331
$ws.Cells.Item(24, 8).Value = 0.879645551787803
$ws.Cells.Item(24, 11).Value = 0.45564956369887
$ws.Cells.Item(24, 13).Value = 0.3213464062991562
$ws.Cells.Item(24, 14).Value = 1.693594944533093

$ws.Cells.Item(25, 2).Value = 0.4204338984738172
$ws.Cells.Item(25, 3).Value = 0.04345664532523585
$ws.Cells.Item(25, 5).Value = 0.1633679991752857
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.8313918579229949
$ws.Cells.Item(25, 8).Value = 0.8839250382530253
$ws.Cells.Item(25, 11).Value = 0.3830911674946265
$ws.Cells.Item(25, 13).Value = 0.2758864611792902
$ws.Cells.Item(25, 14).Value = 1.726903949700041
